# Update header E1 to the "submitted" timestamp, and add a new column F
# (header timestamp in F1, "Submitted" flags in F2/F3) — 3rd version with
# the submit along with date flags.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E1 header text changes from "Status" to a submission timestamp.
$ws.Range("E1").Value = "2025-04-06 15:38:19"

# New column F mirrors column E: header timestamp + per-row "Submitted" flag.
$ws.Range("F1").Value = "2025-04-06 15:39:34"
$ws.Range("F2").Value = "Submitted"
$ws.Range("F3").Value = "Submitted"

# F1 should carry the same header formatting as the rest of row 1 (E1, D1, ...).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats
